$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F2 "Heure" - copy formatting from the adjacent header cell E2
# then set its value (PasteSpecial brings over the bold/centered/bordered style).
$ws.Range("E2").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("F2").Value = "Heure"

# Row 4 (player rank 2)
$ws.Range("B4").Value = "Côme"
$ws.Range("D4").Value = "Sylvie P"
$ws.Range("F4").Value = "13:57"

# Row 5 (player rank 3)
$ws.Range("B5").Value = "Côme"
$ws.Range("D5").Value = "Côme"
$ws.Range("F5").Value = "13:56"

# Row 6 (player rank 4)
$ws.Range("B6").Value = "Côme"
$ws.Range("D6").Value = "Hugo D"
$ws.Range("F6").Value = "13:55"

# Row 7 (player rank 5)
$ws.Range("B7").Value = "Côme"
$ws.Range("D7").Value = "Eric"
$ws.Range("F7").Value = "13:55"

# Row 8 (player rank 6) - D8 changes from "Mathieu" to "Eric", plus new F8
$ws.Range("D8").Value = "Eric"
$ws.Range("F8").Value = "13:55"

# Row 9 (player rank 7) - D9 changes from "Didier" to "Côme", plus new F9
$ws.Range("D9").Value = "Côme"
$ws.Range("F9").Value = "13:55"

# Row 10 (player rank 8) - D10 changes from "Baptiste" to "Côme", plus new F10
$ws.Range("D10").Value = "Côme"
$ws.Range("F10").Value = "13:55"
